$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '63.030.59'
$ws.Cells.Item(2, 5).Value = '  -2.04%  '

$ws.Cells.Item(3, 4).Value = '2.575.26'
$ws.Cells.Item(3, 5).Value = '  -2.97%  '

$c = $ws.Cells.Item(5, 4)
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '587.04'
$c.Style = $s
$ws.Cells.Item(5, 5).Value = '  -3.18%  '

$c = $ws.Cells.Item(6, 4)
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '147.83'
$c.Style = $s
$ws.Cells.Item(6, 5).Value = '  -3.38%  '

$ws.Cells.Item(7, 5).Value = '  +0.02%  '

$ws.Cells.Item(8, 5).Value = '  -1.40%  '

$ws.Cells.Item(9, 5).Value = '  -1.32%  '

$c = $ws.Cells.Item(10, 4)
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '5.71'
$c.Style = $s
$ws.Cells.Item(10, 5).Value = '  +1.63%  '

$c = $ws.Cells.Item(11, 4)
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.380'
$c.Style = $s
$ws.Cells.Item(11, 5).Value = '  -1.95%  '

$ws.Cells.Item(12, 5).Value = '  -0.87%  '

$c = $ws.Cells.Item(13, 4)
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '27.25'
$c.Style = $s
$ws.Cells.Item(13, 5).Value = '  -4.17%  '

$ws.Cells.Item(14, 4).Value = '3.035.42'
$ws.Cells.Item(14, 5).Value = '  -3.04%  '

$ws.Cells.Item(15, 4).Value = '62.903.70'
$ws.Cells.Item(15, 5).Value = '  -1.99%  '

$ws.Cells.Item(16, 5).Value = '  +2.45%  '

$ws.Cells.Item(17, 4).Value = '2.581.78'
$ws.Cells.Item(17, 5).Value = '  -2.31%  '

$c = $ws.Cells.Item(18, 4)
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '12.10'
$c.Style = $s
$ws.Cells.Item(18, 5).Value = '  -0.46%  '

$ws.Cells.Item(19, 5).Value = '  -0.27%  '

$c = $ws.Cells.Item(20, 4)
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '342.88'
$c.Style = $s
$ws.Cells.Item(20, 5).Value = '  -1.45%  '

$ws.Cells.Item(21, 5).Value = '  -2.41%  '

$ws.Cells.Item(22, 5).Value = '  -0.20%  '

$c = $ws.Cells.Item(23, 4)
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '66.55'
$c.Style = $s
$ws.Cells.Item(23, 5).Value = '  -0.21%  '

$c = $ws.Cells.Item(24, 4)
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '1.69'
$c.Style = $s
$ws.Cells.Item(24, 5).Value = '  -3.33%  '

$c = $ws.Cells.Item(25, 4)
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '9.05'
$c.Style = $s
$ws.Cells.Item(25, 5).Value = '  -3.63%  '

$ws.Cells.Item(26, 5).Value = '  -4.16%  '

$c = $ws.Cells.Item(27, 4)
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '553.03'
$c.Style = $s
$ws.Cells.Item(27, 5).Value = '  +0.47%  '

$c = $ws.Cells.Item(28, 4)
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '7.99'
$c.Style = $s
$ws.Cells.Item(28, 5).Value = '  -2.07%  '

$ws.Cells.Item(29, 2).Value = 'Kaspa'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Cells.Item(29, 4)
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.161'
$c.Style = $s
$ws.Cells.Item(29, 5).Value = '  -2.30%  '

$ws.Cells.Item(30, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$c = $ws.Cells.Item(30, 4)
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = $s
$ws.Cells.Item(30, 5).Value = '  +0.17%  '

$c = $ws.Cells.Item(31, 4)
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '2.00'
$c.Style = $s
$ws.Cells.Item(31, 5).Value = '  -2.36%  '

$ws.Cells.Item(32, 5).Value = '  -2.22%  '

$ws.Cells.Item(33, 5).Value = '  -3.15%  '

$c = $ws.Cells.Item(34, 4)
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '5.13'
$c.Style = $s
$ws.Cells.Item(34, 5).Value = '  -4.00%  '

$c = $ws.Cells.Item(35, 4)
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '165.29'
$c.Style = $s
$ws.Cells.Item(35, 5).Value = '  -2.34%  '

$ws.Cells.Item(36, 5).Value = '  -0.08%  '

$c = $ws.Cells.Item(37, 4)
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.998'
$c.Style = $s
$ws.Cells.Item(37, 5).Value = '  -0.17%  '

$c = $ws.Cells.Item(38, 4)
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '19.28'
$c.Style = $s
$ws.Cells.Item(38, 5).Value = '  -0.74%  '

$ws.Cells.Item(39, 5).Value = '  -4.28%  '

$ws.Cells.Item(40, 5).Value = '  +0.02%  '

$c = $ws.Cells.Item(41, 4)
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '164.98'
$c.Style = $s
$ws.Cells.Item(41, 5).Value = '  -0.91%  '

$ws.Cells.Item(42, 2).Value = 'Filecoin'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Cells.Item(42, 4)
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '3.93'
$c.Style = $s
$ws.Cells.Item(42, 5).Value = '  +1.68%  '

$ws.Cells.Item(43, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Cells.Item(43, 4)
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '22.51'
$c.Style = $s
$ws.Cells.Item(43, 5).Value = '  +2.71%  '

$ws.Cells.Item(44, 5).Value = '  +0.51%  '

$ws.Cells.Item(45, 2).Value = 'Mantle'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$c = $ws.Cells.Item(45, 4)
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.626'
$c.Style = $s
$ws.Cells.Item(45, 5).Value = '  -0.88%  '

$ws.Cells.Item(46, 2).Value = 'dogwifhat'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Cells.Item(46, 4)
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '2.02'
$c.Style = $s
$ws.Cells.Item(46, 5).Value = '  +1.11%  '

$ws.Cells.Item(47, 2).Value = 'VeChain'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Cells.Item(47, 4)
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.0246'
$c.Style = $s
$ws.Cells.Item(47, 5).Value = '  -0.84%  '

$ws.Cells.Item(48, 2).Value = 'Stellar'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Cells.Item(48, 4)
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.0958'
$c.Style = $s
$ws.Cells.Item(48, 5).Value = '  -1.03%  '

$ws.Cells.Item(49, 2).Value = 'EnergySwap'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Cells.Item(49, 4)
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '18.83'
$c.Style = $s
$ws.Cells.Item(49, 5).Value = '  -1.32%  '

$ws.Cells.Item(50, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(50, 4).Value = '0.0₆0224'
$ws.Cells.Item(50, 5).Value = '  +10.83%  '

$ws.Cells.Item(51, 2).Value = 'TheGraph'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$c = $ws.Cells.Item(51, 4)
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.178'
$c.Style = $s
$ws.Cells.Item(51, 5).Value = '  -4.51%  '
